# Trade places of the 4th (column F) and 5th (column G) iteration contents
# for the comparison table rows 5-8, and make sure the newly touched
# column I cells pick up the same formatting that is already used by
# the analogous cells further down the sheet (B17 / C17), exactly like
# the source workbook does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- swap the "Iteration IV" / "Iteration V" cell contents -----------------
foreach ($r in 5..8) {
    $colF = $ws.Cells.Item($r, 6)
    $colG = $ws.Cells.Item($r, 7)
    $fVal = $colF.Value()
    $gVal = $colG.Value()
    $colF.Value = $gVal
    $colG.Value = $fVal
}

# --- bring column I formatting in line (mirrors existing styles) -----------
$ws.Range("B17").Copy()
$ws.Range("I5").PasteSpecial(-4122)

$ws.Range("C17").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- restore the cursor/selection position ---------------------------------
[void]$ws.Range("L8").Select()
